# Rerun vaep ranking generator for specific type ids only
#
# This script updates the "R-VAEP" (column C) rankings and Win/Lose (column E)
# labels on the three existing "Top 10 ..." sheets, then adds two brand new
# sheets ("Top 10 VAEP ranking" and "Top 10 proposed ranking") with the
# freshly generated top-10 tables.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Top 10 players goal 90" (sheet2) -- update R-VAEP (C) + one Win/Lose (E)
# ---------------------------------------------------------------------------
$wsGoal = $wb.Worksheets.Item("Top 10 players goal 90")
$wsGoal.Range("C2").Value = 142
$wsGoal.Range("C3").Value = 12
$wsGoal.Range("C4").Value = 28
$wsGoal.Range("C5").Value = 5
$wsGoal.Range("C6").Value = 63
$wsGoal.Range("C7").Value = 33
$wsGoal.Range("E7").Value = "Lose"
$wsGoal.Range("C8").Value = 27
$wsGoal.Range("C9").Value = 17
$wsGoal.Range("C10").Value = 26
$wsGoal.Range("C11").Value = 8

# ---------------------------------------------------------------------------
# 2. "Top 10 players assist 90" (sheet3) -- update R-VAEP (C) + Win/Lose (E)
# ---------------------------------------------------------------------------
$wsAssist = $wb.Worksheets.Item("Top 10 players assist 90")
$wsAssist.Range("C2").Value = 56
$wsAssist.Range("E2").Value = "Win"
$wsAssist.Range("C3").Value = 147
$wsAssist.Range("C4").Value = 5
$wsAssist.Range("C5").Value = 18
$wsAssist.Range("C6").Value = 146
$wsAssist.Range("C7").Value = 139
$wsAssist.Range("E7").Value = "Win"
$wsAssist.Range("C8").Value = 25
$wsAssist.Range("E8").Value = "Lose"
$wsAssist.Range("C9").Value = 33
$wsAssist.Range("E9").Value = "Lose"
$wsAssist.Range("C10").Value = 8
$wsAssist.Range("C11").Value = 23

# ---------------------------------------------------------------------------
# 3. "Top 10 players goal assist 90" (sheet4) -- update R-VAEP (C) + Win/Lose
# ---------------------------------------------------------------------------
$wsGoalAssist = $wb.Worksheets.Item("Top 10 players goal assist 90")
$wsGoalAssist.Range("C2").Value = 56
$wsGoalAssist.Range("E2").Value = "Win"
$wsGoalAssist.Range("C3").Value = 5
$wsGoalAssist.Range("C4").Value = 142
$wsGoalAssist.Range("C5").Value = 33
$wsGoalAssist.Range("C6").Value = 12
$wsGoalAssist.Range("C7").Value = 28
$wsGoalAssist.Range("C8").Value = 147
$wsGoalAssist.Range("C9").Value = 63
$wsGoalAssist.Range("C10").Value = 27
$wsGoalAssist.Range("C11").Value = 18

# ---------------------------------------------------------------------------
# 4. Update selection/scroll state on the three sheets above
# ---------------------------------------------------------------------------
$wsGoal.Activate()
$wsGoal.Range("F17").Select()

$wsAssist.Activate()
$wsAssist.Range("G17").Select()

$wsGoalAssist.Activate()
$wsGoalAssist.Range("D13").Select()

# ---------------------------------------------------------------------------
# 5. New sheet: "Top 10 VAEP ranking"
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsVaep = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsVaep.Name = "Top 10 VAEP ranking"

$wsVaep.Range("A1").Value = "Ranking VAEP"
$wsVaep.Range("B1").Value = "Player Name"
$wsVaep.Range("C1").Value = "Market Value"
$wsVaep.Range("A1:C1").Font.Bold = $true
$wsVaep.Range("A1:C1").HorizontalAlignment = -4108
$wsVaep.Range("A1:C1").WrapText = $true

$vaepRows = @(
    @(1, "Steven Berghuis", "15 million euro"),
    @(2, "Timo Werner", "42 million euro"),
    @(3, "Enis Bardhi", "8 million euro"),
    @(4, "Stuart Armstrong", "6 million euro"),
    @(5, "Adam Szalai", "0,5 million euro"),
    @(6, "Kingsley Coman", "60 million euro"),
    @(7, "Kevin Varga", "2,1 million euro"),
    @(8, "Paul Pogba", "48 million euro"),
    @(9, "Bruno Fernandes", "85 million euro"),
    @(10, "Eric Garcia", "18 million euro")
)
foreach ($row in $vaepRows) {
    $r = $row[0] + 1
    $wsVaep.Cells.Item($r, 1).Value = $row[0]
    $wsVaep.Cells.Item($r, 2).Value = $row[1]
    $wsVaep.Cells.Item($r, 3).Value = $row[2]
}
$wsVaep.Range("C13").Value = "284,6 million euro"

$wsVaep.Columns.Item(1).ColumnWidth = 15.44
$wsVaep.Columns.Item(2).ColumnWidth = 18.26
$wsVaep.Columns.Item(3).ColumnWidth = 17.26

$wsVaep.Activate()
$excel.ActiveWindow.FreezePanes = $true
$wsVaep.Range("C14").Select()

# ---------------------------------------------------------------------------
# 6. New sheet: "Top 10 proposed ranking"
# ---------------------------------------------------------------------------
$wsProposed = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsVaep)
$wsProposed.Name = "Top 10 proposed ranking"

$wsProposed.Range("A1").Value = "Ranking Proposed"
$wsProposed.Range("B1").Value = "Player Name"
$wsProposed.Range("C1").Value = "Market Value"
$wsProposed.Range("A1:C1").Font.Bold = $true
$wsProposed.Range("A1:C1").HorizontalAlignment = -4108
$wsProposed.Range("A1:C1").WrapText = $true

$proposedRows = @(
    @(1, "Patrik Hrosovsky", "5 million euro"),
    @(2, "Laszlo Benes", "8 million euro"),
    @(3, "Ousmane Dembele", "70 million euro"),
    @(4, "Steven Berghuis", "15 million euro"),
    @(5, "Dominic Calvert-Lewin", "30 million euro"),
    @(6, "Dejan Kulusevski", "35 million euro"),
    @(7, "Sergio Busquets", "28 million euro"),
    @(8, "Thiago Alcantara", "48 million euro"),
    @(9, "Cesar Azpilicueta", "24 million euro"),
    @(10, "Sime Vrsaljko", "12 million euro")
)
foreach ($row in $proposedRows) {
    $r = $row[0] + 1
    $wsProposed.Cells.Item($r, 1).Value = $row[0]
    $wsProposed.Cells.Item($r, 2).Value = $row[1]
    $wsProposed.Cells.Item($r, 3).Value = $row[2]
}
$wsProposed.Range("C13").Value = "275 million euro"

$wsProposed.Columns.Item(1).ColumnWidth = 17.71
$wsProposed.Columns.Item(2).ColumnWidth = 19.71
$wsProposed.Columns.Item(3).ColumnWidth = 19.07

$wsProposed.Activate()
$excel.ActiveWindow.FreezePanes = $true
$wsProposed.Range("C14").Select()
